# Update avatar filenames in column H (Avatar) from .jpg to .webp extension
# for the members/pets whose picture asset was re-uploaded in .webp format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "H9"  = "tri.webp"
    "H10" = "khoi.webp"
    "H11" = "phuc.webp"
    "H12" = "thai.webp"
    "H13" = "bong.webp"
    "H14" = "doan.webp"
    "H15" = "quan.webp"
    "H16" = "anh.webp"
    "H17" = "nhu.webp"
    "H18" = "que.webp"
    "H19" = "an.webp"
    "H20" = "mini.webp"
    "H21" = "thanh.webp"
    "H22" = "thuyen.webp"
    "H23" = "bao1.webp"
    "H24" = "tuananh.webp"
    "H25" = "long.webp"
    "H26" = "thien.webp"
    "H27" = "hiep.webp"
    "H38" = "honganh.webp"
    "H41" = "bedau.webp"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# Reflect the cell that was active/selected when the workbook was last saved.
$ws.Range("H38").Select()
